$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '63.625.72'
Set-TextCell $ws.Range('E2') '  +3.09%  '
Set-TextCell $ws.Range('D3') '3.484.84'
Set-TextCell $ws.Range('E3') '  +2.15%  '
Set-TextCell $ws.Range('E4') '  +0.16%  '
Set-TextCell $ws.Range('D5') '582.17'
Set-TextCell $ws.Range('E5') '  +0.89%  '
Set-TextCell $ws.Range('D6') '147.64'
Set-TextCell $ws.Range('E6') '  +2.29%  '
Set-TextCell $ws.Range('D7') '3.480.18'
Set-TextCell $ws.Range('E7') '  +2.01%  '
Set-TextCell $ws.Range('E8') '  +0.07%  '
Set-TextCell $ws.Range('D9') '0.478'
Set-TextCell $ws.Range('E9') '  +0.96%  '
Set-TextCell $ws.Range('E11') '  +2.10%  '
Set-TextCell $ws.Range('D12') '0.406'
Set-TextCell $ws.Range('E12') '  +5.76%  '
Set-TextCell $ws.Range('D13') '4.088.68'
Set-TextCell $ws.Range('E13') '  +2.37%  '
Set-TextCell $ws.Range('D14') '29.86'
Set-TextCell $ws.Range('E14') '  +6.82%  '
Set-TextCell $ws.Range('D15') '0.128'
Set-TextCell $ws.Range('E15') '  +2.65%  '
Set-TextCell $ws.Range('D16') '3.480.28'
Set-TextCell $ws.Range('E16') '  +2.03%  '
Set-TextCell $ws.Range('D17') '0.0000172'
Set-TextCell $ws.Range('E17') '  +1.83%  '
Set-TextCell $ws.Range('D18') '63.673.31'
Set-TextCell $ws.Range('E18') '  +3.15%  '
Set-TextCell $ws.Range('D19') '6.38'
Set-TextCell $ws.Range('E19') '  +4.12%  '
Set-TextCell $ws.Range('D20') '14.47'
Set-TextCell $ws.Range('E20') '  +4.70%  '
Set-TextCell $ws.Range('D21') '9.40'
Set-TextCell $ws.Range('E21') '  +2.91%  '
Set-TextCell $ws.Range('D22') '391.06'
Set-TextCell $ws.Range('E22') '  +0.54%  '
Set-TextCell $ws.Range('E23') '  +3.30%  '
Set-TextCell $ws.Range('D24') '75.27'
Set-TextCell $ws.Range('E24') '  +1.44%  '
Set-TextCell $ws.Range('E25') '  -0.11%  '
Set-TextCell $ws.Range('D26') '3.623.32'
Set-TextCell $ws.Range('E26') '  +2.05%  '
Set-TextCell $ws.Range('E27') '  +2.09%  '
Set-TextCell $ws.Range('D28') '0.180'
Set-TextCell $ws.Range('E28') '  -4.63%  '
Set-TextCell $ws.Range('E29') '  +3.66%  '
Set-TextCell $ws.Range('E30') '  +0.16%  '
Set-TextCell $ws.Range('D31') '8.26'
Set-TextCell $ws.Range('E31') '  +3.52%  '
Set-TextCell $ws.Range('E32') '  +0.09%  '
Set-TextCell $ws.Range('E33') '  +0.08%  '
Set-TextCell $ws.Range('E34') '  -2.13%  '
Set-TextCell $ws.Range('D35') '23.67'
Set-TextCell $ws.Range('E35') '  +1.00%  '
Set-TextCell $ws.Range('D36') '7.16'
Set-TextCell $ws.Range('E36') '  +3.23%  '
Set-TextCell $ws.Range('D37') '5.32'
Set-TextCell $ws.Range('E37') '  +2.86%  '
Set-TextCell $ws.Range('B38') 'ImmutableX'
Set-TextCell $ws.Range('C38') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws.Range('D38') '1.59'
Set-TextCell $ws.Range('E38') '  +8.34%  '
Set-TextCell $ws.Range('B39') 'EnergySwap'
Set-TextCell $ws.Range('C39') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range('D39') '31.62'
Set-TextCell $ws.Range('E39') '  +11.08%  '
Set-TextCell $ws.Range('D40') '169.81'
Set-TextCell $ws.Range('E40') '  +0.79%  '
Set-TextCell $ws.Range('D41') '3.528.20'
Set-TextCell $ws.Range('E41') '  +2.46%  '
Set-TextCell $ws.Range('D42') '0.0769'
Set-TextCell $ws.Range('E42') '  +2.08%  '
Set-TextCell $ws.Range('D43') '0.799'
Set-TextCell $ws.Range('E43') '  +1.77%  '
Set-TextCell $ws.Range('D44') '1.75'
Set-TextCell $ws.Range('E44') '  +4.75%  '
Set-TextCell $ws.Range('D45') '42.37'
Set-TextCell $ws.Range('E45') '  -0.43%  '
Set-TextCell $ws.Range('E46') '  +4.17%  '
Set-TextCell $ws.Range('D47') '4.44'
Set-TextCell $ws.Range('E47') '  +0.02%  '
Set-TextCell $ws.Range('D48') '2.629.30'
Set-TextCell $ws.Range('E48') '  +5.08%  '
Set-TextCell $ws.Range('E49') '  +11.57%  '
Set-TextCell $ws.Range('D50') '23.25'
Set-TextCell $ws.Range('E50') '  +2.50%  '
Set-TextCell $ws.Range('D51') '6.80'
Set-TextCell $ws.Range('E51') '  +3.10%  '
